$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.867.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.92%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.66"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.29%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.61%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3691"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.80%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07171"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8774"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.85%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07877"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.44%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "19.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.859.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.32%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.329"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.386"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.90"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.34%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008731"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.24%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.911.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.91%  "

$ws.Range("E21").Value = "  -2.62%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.993"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.89%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.44"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.976"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "150.74"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.49%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.969"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.16%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.49"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.934"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.87%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08823"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.128"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.45%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7554"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.33%  "

$ws.Range("E33").Value = "  -0.66%  "

$ws.Range("E34").Value = "  -2.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.603"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.094"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01937"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.924"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.75%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05127"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.32%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.916"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4979"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.79%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1596"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.10%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.361"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.07%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4671"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.56%  "

$ws.Range("E45").Value = "  +0.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.21%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("E48").Value = "  -2.54%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06096"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "36.44"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.84%  "
